# Update countries & provincias Spain
# Applies the COVID data refresh: numeric updates for a set of rows, four
# country rank swaps (adjacent rows exchanging names because the newer
# snapshot re-sorted them), and a refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 18:49"

# --- Country name swaps (adjacent rows trade rank) ---
$ws.Range("A20").Value = "Irak"
$ws.Range("A21").Value = "Pakistan"

$ws.Range("A89").Value = "Grecia"
$ws.Range("A90").Value = "Croacia"

$ws.Range("A137").Value = "Trinidad yTobago"
$ws.Range("A138").Value = "Aruba"

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Numeric data refresh (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6844149
$ws.Range("C4").Value = 15848
$ws.Range("D4").Value = 4133363
$ws.Range("E4").Value = 2509087
$ws.Range("G4").Value = 351
$ws.Range("H4").Value = 201699

# Row 5 - India
$ws.Range("B5").Value = 5169981
$ws.Range("C5").Value = 54088
$ws.Range("D5").Value = 4068260
$ws.Range("E5").Value = 1017986
$ws.Range("G5").Value = 505
$ws.Range("H5").Value = 83735

# Row 6 - Brasil
$ws.Range("B6").Value = 4430227
$ws.Range("C6").Value = 8541
$ws.Range("E6").Value = 575552
$ws.Range("G6").Value = 189
$ws.Range("H6").Value = 134363

# Row 20 - Irak (post-swap)
$ws.Range("B20").Value = 307385
$ws.Range("C20").Value = 4326
$ws.Range("D20").Value = 241100
$ws.Range("E20").Value = 57953
$ws.Range("G20").Value = 84
$ws.Range("H20").Value = 8332

# Row 21 - Pakistan (post-swap)
$ws.Range("B21").Value = 303634
$ws.Range("C21").Value = 545
$ws.Range("D21").Value = 291169
$ws.Range("E21").Value = 6066
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 6399

# Row 23 - Filipinas
$ws.Range("B23").Value = 293025
$ws.Range("C23").Value = 1585
$ws.Range("D23").Value = 215954
$ws.Range("E23").Value = 41413
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 35658

# Row 29 - Israel
$ws.Range("B29").Value = 140539
$ws.Range("C29").Value = 792
$ws.Range("D29").Value = 122836
$ws.Range("E29").Value = 8504
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 9199

# Row 32 - Rumania
$ws.Range("B32").Value = 122257
$ws.Range("C32").Value = 732
$ws.Range("E32").Value = 14165
$ws.Range("G32").Value = 33
$ws.Range("H32").Value = 11029

# Row 66 - Azerbaiyan
$ws.Range("B66").Value = 42739
$ws.Range("C66").Value = 1707
$ws.Range("D66").Value = 23321
$ws.Range("E66").Value = 18930
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 488

# Row 89 - Grecia (post-swap)
$ws.Range("B89").Value = 14400
$ws.Range("C89").Value = 359
$ws.Range("D89").Value = 3804
$ws.Range("E89").Value = 10271
$ws.Range("G89").Value = 9
$ws.Range("H89").Value = 325

# Row 90 - Croacia (post-swap)
$ws.Range("B90").Value = 14279
$ws.Range("C90").Value = 250
$ws.Range("D90").Value = 11933
$ws.Range("E90").Value = 2108
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 238

# Row 93 - Republica de Macedonia
$ws.Range("B93").Value = 12571
$ws.Range("C93").Value = 73
$ws.Range("E93").Value = 1934
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 266

# Row 94 - Madagascar
$ws.Range("B94").Value = 11948
$ws.Range("C94").Value = 132
$ws.Range("D94").Value = 6788
$ws.Range("E94").Value = 4813
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 347

# Row 137 - Trinidad yTobago (post-swap)
$ws.Range("B137").Value = 3336
$ws.Range("C137").Value = 9
$ws.Range("D137").Value = 814
$ws.Range("E137").Value = 2464
$ws.Range("H137").Value = 58

# Row 138 - Aruba (post-swap)
$ws.Range("B138").Value = 3328
$ws.Range("D138").Value = 1676
$ws.Range("E138").Value = 1630
$ws.Range("H138").Value = 22

# Row 162 - Republica de Chipre
$ws.Range("B162").Value = 1333
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 37

# Row 195 - Monaco
$ws.Range("B195").Value = 112
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 108
$ws.Range("E195").Value = 3

# Row 214 - Islas Malvinas (post-swap)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215 - Montserrat (post-swap)
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
